$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1369.6
$ws.Range("I70").Value = 1024.5
$ws.Range("K70").Value = 3073.5
$ws.Range("M70").Value = -2803.5
$ws.Range("H73").Value = 1369.6
$ws.Range("I73").Value = 1024.5
$ws.Range("K73").Value = 3073.5
$ws.Range("M73").Value = -2137.5
$ws.Range("H86").Value = 1654.2858
$ws.Range("I86").Value = 1654.2858
$ws.Range("K86").Value = 1654.2858
$ws.Range("M86").Value = -531.2858000000001
$ws.Range("H89").Value = 1654.2858
$ws.Range("I89").Value = 1654.2858
$ws.Range("K89").Value = 8271.429
$ws.Range("M89").Value = -2655.429
$ws.Range("H92").Value = 570.75
$ws.Range("I92").Value = 427.83334
$ws.Range("K92").Value = 427.83334
$ws.Range("M92").Value = 820.16666
$ws.Range("H111").Value = 997.5263
$ws.Range("I111").Value = 934.0833
$ws.Range("J111").Value = 1106.2858
$ws.Range("K111").Value = 2802.2499
$ws.Range("L111").Value = 3318.8574
$ws.Range("M111").Value = 264.7501000000002
$ws.Range("N111").Value = -9452.857400000001
$ws.Range("H132").Value = 4905.8184
$ws.Range("I132").Value = 4896
$ws.Range("K132").Value = 14688
$ws.Range("M132").Value = -12158
$ws.Range("H140").Value = 74999.5
$ws.Range("J140").Value = 74999.5
$ws.Range("L140").Value = 74999.5
$ws.Range("N140").Value = -85359.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1983.6
$ws.Range("J45").Value = 2175
$ws.Range("L45").Value = 2175
$ws.Range("N45").Value = -2929
$ws.Range("H97").Value = 735.7778
$ws.Range("I97").Value = 1599.6666
$ws.Range("J97").Value = 303.83334
$ws.Range("K97").Value = 1599.6666
$ws.Range("L97").Value = 303.83334
$ws.Range("M97").Value = -1103.6666
$ws.Range("N97").Value = -1295.83334
$ws.Range("H140").Value = 72151.2
$ws.Range("J140").Value = 71439
$ws.Range("L140").Value = 71439
$ws.Range("N140").Value = -81799
$ws.Range("H141").Value = 48962.2
$ws.Range("I141").Value = 51219
$ws.Range("J141").Value = 47457.668
$ws.Range("K141").Value = 51219
$ws.Range("L141").Value = 47457.668
$ws.Range("M141").Value = -46039
$ws.Range("N141").Value = -57817.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 417
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 417
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 417
$ws.Range("N94").Value = -1319
$ws.Range("H99").Value = 4519
$ws.Range("I99").Value = 4378
$ws.Range("K99").Value = 4378
$ws.Range("M99").Value = -2880
$ws.Range("H139").Value = 74000
$ws.Range("J139").Value = 74000
$ws.Range("L139").Value = 74000
$ws.Range("N139").Value = -84280
$ws.Range("H140").Value = 69000
$ws.Range("J140").Value = 69000
$ws.Range("L140").Value = 69000
$ws.Range("N140").Value = -79360

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3146.1875
$ws.Range("I22").Value = 1593.6
$ws.Range("J22").Value = 3851.9092
$ws.Range("K22").Value = 1593.6
$ws.Range("L22").Value = 3851.9092
$ws.Range("M22").Value = -1243.6
$ws.Range("N22").Value = -4551.9092
$ws.Range("H58").Value = 3398.5
$ws.Range("I58").Value = 2900
$ws.Range("K58").Value = 2900
$ws.Range("M58").Value = -2697
$ws.Range("H62").Value = 13624.542
$ws.Range("I62").Value = 9587.571
$ws.Range("K62").Value = 9587.571
$ws.Range("M62").Value = -8963.571
$ws.Range("H65").Value = 13624.542
$ws.Range("I65").Value = 9587.571
$ws.Range("K65").Value = 47937.855
$ws.Range("M65").Value = -44817.855
$ws.Range("H105").Value = 2291.8823
$ws.Range("I105").Value = 1783.5333
$ws.Range("K105").Value = 1783.5333
$ws.Range("M105").Value = -36.53330000000005
$ws.Range("H122").Value = 1170
$ws.Range("I122").Value = 853.5
$ws.Range("K122").Value = 2560.5
$ws.Range("M122").Value = -110.5
$ws.Range("H132").Value = 1893.8182
$ws.Range("I132").Value = 1872.75
$ws.Range("K132").Value = 5618.25
$ws.Range("M132").Value = -3088.25
$ws.Range("H136").Value = 3398.5
$ws.Range("I136").Value = 2900
$ws.Range("K136").Value = 8700
$ws.Range("M136").Value = -6150

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 723.3
$ws.Range("I8").Value = 723.3
$ws.Range("K8").Value = 2169.9
$ws.Range("M8").Value = -2030.9
$ws.Range("H57").Value = 21000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H68").Value = 1170.625
$ws.Range("I68").Value = 1199.5714
$ws.Range("J68").Value = 968
$ws.Range("K68").Value = 3598.7142
$ws.Range("L68").Value = 2904
$ws.Range("M68").Value = -2787.7142
$ws.Range("N68").Value = -4526
$ws.Range("H71").Value = 1170.625
$ws.Range("I71").Value = 1199.5714
$ws.Range("J71").Value = 968
$ws.Range("K71").Value = 10796.1426
$ws.Range("L71").Value = 8712
$ws.Range("M71").Value = -6740.142600000001
$ws.Range("N71").Value = -16824
$ws.Range("H102").Value = 4837.875
$ws.Range("H105").Value = 36962.668
$ws.Range("J105").Value = 36962.668
$ws.Range("L105").Value = 110888.004
$ws.Range("N105").Value = -116130.004
$ws.Range("H108").Value = 2666.5
$ws.Range("I108").Value = 2666.5
$ws.Range("K108").Value = 7999.5
$ws.Range("M108").Value = -5119.5
$ws.Range("H109").Value = 3647.4285
$ws.Range("I109").Value = 3922
$ws.Range("K109").Value = 11766
$ws.Range("M109").Value = -10726
$ws.Range("H111").Value = 9999
$ws.Range("I111").Value = 9999
$ws.Range("K111").Value = 29997
$ws.Range("M111").Value = -26930
$ws.Range("H112").Value = 4257.5
$ws.Range("I112").Value = 1000
$ws.Range("K112").Value = 3000
$ws.Range("M112").Value = -1892
$ws.Range("H113").Value = 678.5714
$ws.Range("J113").Value = 874.25
$ws.Range("L113").Value = 2622.75
$ws.Range("N113").Value = -6962.75
$ws.Range("H114").Value = 3206.6667
$ws.Range("I114").Value = 474
$ws.Range("J114").Value = 4573
$ws.Range("K114").Value = 1422
$ws.Range("L114").Value = 13719
$ws.Range("M114").Value = 1832
$ws.Range("N114").Value = -20227
$ws.Range("H115").Value = 2902.1428
$ws.Range("I115").Value = 2453.75
$ws.Range("K115").Value = 7361.25
$ws.Range("M115").Value = -6186.25
$ws.Range("H122").Value = 38324.57
$ws.Range("J122").Value = 41195.69
$ws.Range("L122").Value = 370761.21
$ws.Range("N122").Value = -375661.21

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1289.5
$ws.Range("I97").Value = 1289.5
$ws.Range("K97").Value = 1289.5
$ws.Range("M97").Value = -793.5
$ws.Range("H107").Value = 625.4211
$ws.Range("I107").Value = 281.7143
$ws.Range("K107").Value = 281.7143
$ws.Range("M107").Value = 1638.2857
$ws.Range("H140").Value = 69770
$ws.Range("J140").Value = 69770
$ws.Range("L140").Value = 69770
$ws.Range("N140").Value = -80130

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2691.5
$ws.Range("I10").Value = 2941
$ws.Range("J10").Value = 2566.75
$ws.Range("K10").Value = 2941
$ws.Range("L10").Value = 2566.75
$ws.Range("M10").Value = -2801
$ws.Range("N10").Value = -2846.75
$ws.Range("H22").Value = 469.7143
$ws.Range("I22").Value = 438
$ws.Range("J22").Value = 549
$ws.Range("K22").Value = 438
$ws.Range("L22").Value = 549
$ws.Range("M22").Value = -143
$ws.Range("N22").Value = -1139
$ws.Range("H27").Value = 469.7143
$ws.Range("I27").Value = 438
$ws.Range("J27").Value = 549
$ws.Range("K27").Value = 438
$ws.Range("L27").Value = 549
$ws.Range("M27").Value = -331
$ws.Range("N27").Value = -763
$ws.Range("H132").Value = 2839.6
$ws.Range("I132").Value = 2299.5
$ws.Range("K132").Value = 6898.5
$ws.Range("M132").Value = -4368.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 13411.5
$ws.Range("I37").Value = 16633
$ws.Range("J37").Value = 12030.857
$ws.Range("K37").Value = 16633
$ws.Range("L37").Value = 12030.857
$ws.Range("M37").Value = -16430
$ws.Range("N37").Value = -12436.857
$ws.Range("H70").Value = 27929.578
$ws.Range("I70").Value = 15900
$ws.Range("K70").Value = 15900
$ws.Range("M70").Value = -15585
$ws.Range("H73").Value = 27929.578
$ws.Range("I73").Value = 15900
$ws.Range("K73").Value = 15900
$ws.Range("M73").Value = -14808
$ws.Range("H96").Value = 3057.5715
$ws.Range("I96").Value = 10003
$ws.Range("K96").Value = 10003
$ws.Range("M96").Value = -8630
$ws.Range("H122").Value = 5426.2104
$ws.Range("J122").Value = 4652.5
$ws.Range("L122").Value = 13957.5
$ws.Range("N122").Value = -18857.5
$ws.Range("H126").Value = 2974.125
$ws.Range("I126").Value = 1368
$ws.Range("K126").Value = 4104
$ws.Range("M126").Value = -1634
$ws.Range("H132").Value = 1959.7778
$ws.Range("I132").Value = 1959.7778
$ws.Range("K132").Value = 5879.3334
$ws.Range("M132").Value = -3349.3334

# ---- Cell deletions (values removed entirely) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M94").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M57").ClearContents()
